$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.548.85'
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").Value = '1.876.49'
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.015'
$ws.Range("E4").Value = '  +1.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.74'
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4791'
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3781'
$ws.Range("E8").Value = '  +2.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07384'
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9392'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.78'
$ws.Range("E11").Value = '  +5.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07862'
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("D13").Value = '1.890.63'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.452'
$ws.Range("E14").Value = '  +2.59%  '
$ws.Range("E15").Value = '  +3.02%  '
$ws.Range("E16").Value = '  +2.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.016'
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008965'
$ws.Range("E18").Value = '  +3.58%  '
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.95'
$ws.Range("E20").Value = '  +2.78%  '
$ws.Range("D21").Value = '27.599.34'
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.142'
$ws.Range("E22").Value = '  +2.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.76'
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.957'
$ws.Range("E24").Value = '  +2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.78'
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.58'
$ws.Range("E26").Value = '  +2.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.021'
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.04'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.021'
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08941'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.327'
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("E32").Value = '  +3.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.613'
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7509'
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.694'
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02075'
$ws.Range("E36").Value = '  +6.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.118'
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05309'
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.010'
$ws.Range("E39").Value = '  +1.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5358'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.114'
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1526'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.437'
$ws.Range("E43").Value = '  +2.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.65'
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4846'
$ws.Range("E45").Value = '  +3.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.016'
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.662'
$ws.Range("E47").Value = '  +3.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.16'
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.34'
$ws.Range("E49").Value = '  +3.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06102'
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9007'
$ws.Range("E51").Value = '  +1.79%  '
